$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "14.1.0" release column (O) - header + Tardigrade.Framework version bump.
$ws.Range("O1").Value2 = "14.1.0"
$ws.Range("O1").Font.Bold = $true
$ws.Range("O2").Value2 = "12.1.0"

# Fill in the previously-blank MailKit row (row 9) with "-" placeholders for
# the releases where MailKit wasn't applicable yet.
$ws.Range("B9").Value2 = "-"
$ws.Range("C9").Value2 = "-"
$ws.Range("D9").Value2 = "-"
$ws.Range("E9").Value2 = "-"
$ws.Range("F9").Value2 = "-"
$ws.Range("G9").Value2 = "-"
$ws.Range("H9").Value2 = "-"
$ws.Range("I9").Value2 = "-"
$ws.Range("J9").Value2 = "-"
$ws.Range("K9").Value2 = "-"
$ws.Range("L9").Value2 = "-"
$ws.Range("M9").Value2 = "-"

# Match the author's final selection/cursor position.
$ws.Range("N9").Select()
